$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06965798242991458
$ws.Range("C2").Value = 0.05838643745354499
$ws.Range("D2").Value = 0.02470440429076765
$ws.Range("E2").Value = 0.02238195385626012
$ws.Range("F2").Value = 0.03729000840050865
$ws.Range("B3").Value = 0.06662487478125043
$ws.Range("C3").Value = 0.05892047593726649
$ws.Range("D3").Value = 0.02667125851262939
$ws.Range("E3").Value = 0.02200709936861224
$ws.Range("F3").Value = 0.03600565181035668
$ws.Range("B4").Value = 0.05563203081583885
$ws.Range("C4").Value = 0.05525538879418047
$ws.Range("D4").Value = 0.03029388678197968
$ws.Range("E4").Value = 0.02973899039888512
$ws.Range("F4").Value = 0.03781497346398539
$ws.Range("B5").Value = 0.08774816749819409
$ws.Range("C5").Value = 0.06241552089222512
$ws.Range("D5").Value = 0.03662967293334959
$ws.Range("E5").Value = 0.03799231118257523
$ws.Range("F5").Value = 0.04944292076152362
$ws.Range("B6").Value = 0.0883993775006465
$ws.Range("C6").Value = 0.08548604603941656
$ws.Range("D6").Value = 0.03951016831172372
$ws.Range("E6").Value = 0.04483217888356641
$ws.Range("F6").Value = 0.05465443224694995
$ws.Range("B7").Value = 0.07578485615025315
$ws.Range("C7").Value = 0.09814763071669683
$ws.Range("D7").Value = 0.03986108376387477
$ws.Range("E7").Value = 0.03304782072504283
$ws.Range("F7").Value = 0.0588356924477206
$ws.Range("B8").Value = 0.08145110129797478
$ws.Range("C8").Value = 0.0757279231707486
$ws.Range("D8").Value = 0.02868446819530372
$ws.Range("E8").Value = 0.02544855481092779
$ws.Range("F8").Value = 0.0492978030126445
$ws.Range("B9").Value = 0.06560115852706098
$ws.Range("C9").Value = 0.06292960986673535
$ws.Range("D9").Value = 0.02471937779873417
$ws.Range("E9").Value = 0.02234652618971443
$ws.Range("F9").Value = 0.03567285397568205
$ws.Range("B10").Value = 0.07828183596336742
$ws.Range("C10").Value = 0.05682810712890048
$ws.Range("D10").Value = 0.02468435939887132
$ws.Range("E10").Value = 0.01595103324473795
$ws.Range("F10").Value = 0.03765077912020958
$ws.Range("B11").Value = 0.08926538937592576
$ws.Range("C11").Value = 0.05692797099370564
$ws.Range("D11").Value = 0.02705983169499722
$ws.Range("E11").Value = 0.01391552377667354
$ws.Range("F11").Value = 0.03545974688618096
$ws.Range("B12").Value = 0.08971346413468463
$ws.Range("C12").Value = 0.05560678120822612
$ws.Range("D12").Value = 0.02893944091974564
$ws.Range("E12").Value = 0.01558929248172013
$ws.Range("F12").Value = 0.0368747381765368
$ws.Range("B13").Value = 0.09478855953145492
$ws.Range("C13").Value = 0.05672376912003907
$ws.Range("D13").Value = 0.03007283321842252
$ws.Range("E13").Value = 0.01727270456705887
$ws.Range("F13").Value = 0.03755896202336916
$ws.Range("B14").Value = 0.08853842659419367
$ws.Range("C14").Value = 0.05761656478208125
$ws.Range("D14").Value = 0.03156857817738345
$ws.Range("E14").Value = 0.01844083748753854
$ws.Range("F14").Value = 0.03745356065612188
$ws.Range("B15").Value = 0.09259772046836198
$ws.Range("C15").Value = 0.05755085259275139
$ws.Range("D15").Value = 0.03270614913341009
$ws.Range("E15").Value = 0.01809749343429279
$ws.Range("F15").Value = 0.03738662047211606
$ws.Range("B16").Value = 0.09284152478025973
$ws.Range("C16").Value = 0.05737536453395244
$ws.Range("D16").Value = 0.03282332344044964
$ws.Range("E16").Value = 0.0191646522754241
$ws.Range("F16").Value = 0.03765938270519531
$ws.Range("B17").Value = 0.092479601758826
$ws.Range("C17").Value = 0.05744103448055111
$ws.Range("D17").Value = 0.03314485927196179
$ws.Range("E17").Value = 0.01866403548609182
$ws.Range("F17").Value = 0.03785601096323388
$ws.Range("B18").Value = 0.09235694632284203
$ws.Range("C18").Value = 0.05743081157306191
$ws.Range("D18").Value = 0.03315419945604126
$ws.Range("E18").Value = 0.01872433524452861
$ws.Range("F18").Value = 0.03782205485909784
